$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old hyperlink on B2 (mailto to the old regd-number address)
$ws.Range("B2").Hyperlinks.Delete()

# --- Row 1 (headers): keep "name", change regd/Email/Campus/School/Dept
#     header row down to just name / email / campus
$ws.Cells.Item(1, 1).Value = "name"
$ws.Cells.Item(1, 2).Value = "email"
$ws.Cells.Item(1, 3).Value = "campus"
$ws.Range("D1:F1").ClearContents()

# --- Row 2 (data): single remaining person's info
$ws.Cells.Item(2, 1).Value = "chinmaya"
$ws.Cells.Item(2, 2).Value = "situ@chinmayakumarbiswal.in"
$ws.Cells.Item(2, 3).Value = "BBSR"
$ws.Range("D2:F2").ClearContents()

# Re-add the hyperlink on B2 pointing at the new email address, then restore
# the cell's Hyperlink style (Add() bumps it to a fresh, equivalent style)
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:situ@chinmayakumarbiswal.in")
$ws.Range("B2").Style = "Hyperlink"

# --- Column widths for the now-relevant A:C columns
$ws.Columns.Item(1).ColumnWidth = 9.109375
$ws.Columns.Item(2).ColumnWidth = 27.88671875
$ws.Columns.Item(3).ColumnWidth = 7.109375

# --- Selection moves to G9
[void]$ws.Range("G9").Select()
